$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 224; existing rows 224:239 shift down to 225:240
$ws.Rows(224).Insert()

# Populate the new row 224 with the new price-report entry
$ws.Range("A224").Value = 11
$ws.Range("B224").Value = "Vega Monumental Concepción"
$ws.Range("C224").Value = "Bíobío"
$ws.Range("D224").Value = 44931
$ws.Range("E224").Value = 8
$ws.Range("F224").Value = 100112003
$ws.Range("G224").Value = "Ajo"
$ws.Range("H224").Value = "Chino"
$ws.Range("I224").Value = "Primera"
$ws.Range("J224").Value = 220
$ws.Range("K224").Value = 11000
$ws.Range("L224").Value = 12000
$ws.Range("M224").Value = 11545
$ws.Range("N224").Value = "$/caja 10 kilos"
$ws.Range("O224").Value = "China"
$ws.Range("P224").Value = 1154
$ws.Range("Q224").Value = 10
$ws.Range("R224").Value = "Hortaliza"
